$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended below the existing table (row 9): "helix" / "com.singleton.helix"
$ws.Range("A9").Value = "helix"
$ws.Range("B9").Value = "com.singleton.helix"

# Match formatting (style) of the row above it (row 8) by copying formats only
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)

# Move the active selection down to the next empty row, as in the source workbook
$ws.Range("A10").Select()
